$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting so values
# like "99.10" or "  +0.42%  " are not coerced into numbers and
# lose trailing zeros / padding.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '37.751.10'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '2.084.95'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '233.56'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '58.95'
$ws.Range('E7').Value = '  +2.97%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +1.91%  '
$ws.Range('E10').Value = '  +0.20%  '
$ws.Range('E11').Value = '  +2.86%  '
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('E13').Value = '  +2.19%  '
$ws.Range('D14').Value = '21.25'
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('E15').Value = '  +1.87%  '
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('D17').Value = '2.082.93'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').Value = '37.726.20'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('D19').Value = '6.16'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').Value = '71.84'
$ws.Range('E20').Value = '  +1.77%  '
$ws.Range('E21').Value = '  +3.23%  '
$ws.Range('D22').Value = '228.41'
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').Value = '2.41'
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('D26').Value = '9.62'
$ws.Range('E26').Value = '  +7.93%  '
$ws.Range('D27').Value = '171.17'
$ws.Range('E27').Value = '  +0.55%  '
$ws.Range('D28').Value = '0.136'
$ws.Range('E28').Value = '  -1.46%  '
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('D30').Value = '19.58'
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('E31').Value = '  +2.28%  '
$ws.Range('D32').Value = '4.75'
$ws.Range('E32').Value = '  +2.32%  '
$ws.Range('E33').Value = '  +1.68%  '
$ws.Range('D34').Value = '4.68'
$ws.Range('E34').Value = '  +1.26%  '
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('E36').Value = '  +1.08%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = '5.43'
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('E40').Value = '  -1.39%  '
$ws.Range('D41').Value = '17.43'
$ws.Range('E41').Value = '  +11.64%  '
$ws.Range('D42').Value = '99.10'
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('E43').Value = '  +2.64%  '
$ws.Range('D44').Value = '2.91'
$ws.Range('E44').Value = '  -1.01%  '
$ws.Range('D45').Value = '1.452.49'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range('D47').Value = '4.15'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('D49').Value = '7.37'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('D51').Value = '2.275.90'
$ws.Range('E51').Value = '  +0.42%  '
